$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume strings are kept as text (matches source formatting)
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    'D2' = '26.037.75'
    'E2' = '  +2.94%  '
    'D3' = '1.596.60'
    'E3' = '  +2.07%  '
    'E4' = '  -0.12%  '
    'D5' = '212.67'
    'E5' = '  +2.78%  '
    'E6' = '  -0.15%  '
    'E7' = '  +1.53%  '
    'E8' = '  +2.70%  '
    'E9' = '  +1.58%  '
    'D10' = '17.94'
    'E10' = '  +1.11%  '
    'D11' = '0.0816'
    'E11' = '  +4.40%  '
    'D12' = '1.819.45'
    'E12' = '  +2.09%  '
    'D13' = '1.594.25'
    'E13' = '  +1.89%  '
    'E14' = '  -0.28%  '
    'D15' = '0.510'
    'E15' = '  +1.05%  '
    'D16' = '26.024.50'
    'E16' = '  +2.90%  '
    'D17' = '60.34'
    'E17' = '  +1.87%  '
    'E18' = '  +1.43%  '
    'E19' = '  -0.07%  '
    'D20' = '203.21'
    'E20' = '  +9.86%  '
    'E21' = '  +2.68%  '
    'D22' = '9.29'
    'E22' = '  +0.34%  '
    'E23' = '  +1.70%  '
    'E24' = '  +11.45%  '
    'D25' = '140.87'
    'E25' = '  +0.74%  '
    'E26' = '  -0.16%  '
    'D27' = '0.124'
    'E27' = '  -2.65%  '
    'D28' = '15.20'
    'E28' = '  +2.71%  '
    'E29' = '  -0.45%  '
    'E30' = '  +1.49%  '
    'E31' = '  +1.37%  '
    'E32' = '  +3.02%  '
    'E33' = '  -1.24%  '
    'E34' = '  +1.26%  '
    'E35' = '  +1.96%  '
    'D36' = '1.107.10'
    'E36' = '  +1.75%  '
    'E37' = '  +8.27%  '
    'E38' = '  +0.20%  '
    'E39' = '  +0.33%  '
    'E40' = '  +1.99%  '
    'D41' = '0.492'
    'E41' = '  -0.07%  '
    'D42' = '0.777'
    'E42' = '  -2.50%  '
    'D43' = '1.732.17'
    'E43' = '  +2.05%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D44' = '92.59'
    'E44' = '  -0.51%  '
    'B45' = 'FraxShare'
    'C45' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D45' = '5.09'
    'E45' = '  +0.86%  '
    'D46' = '1.50'
    'E46' = '  +4.49%  '
    'D47' = '53.31'
    'E47' = '  +1.70%  '
    'D48' = '0.0503'
    'E48' = '  -0.32%  '
    'E49' = '  +0.68%  '
    'E50' = '  +0.03%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D51' = '7.22'
    'E51' = '  +0.83%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

